$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at positions 8-9, shifting the current extr1..extr8
# rows (originally rows 8-15) down to rows 10-17. The new rows will hold
# "line7" and "line8" entries.
$ws.Rows("8:9").Insert()

# Match the header-style formatting (bold, centered, bordered) used by the
# other index cells in column A for the two freshly inserted rows.
$styleRange = $ws.Range("A8:A9")
$styleRange.Font.Bold = $true
$styleRange.HorizontalAlignment = -4108
$styleRange.VerticalAlignment = -4160
$styleRange.Borders.LineStyle = 1

# Row 8: line7
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# Row 9: line8
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $false

# The rows that used to be 8 (extr1) and 9 (extr2) are now at 10 and 11;
# their in_service flag flips from FALSE to TRUE.
$ws.Cells.Item(10, 5).Value = $true
$ws.Cells.Item(11, 5).Value = $true

# Column A is a simple running 0-based index; renumber rows 10-17 (the
# rows pushed down by the insert) so the sequence keeps counting up
# (8, 9, 10, ... 15) instead of repeating the old 6, 7 values.
for ($r = 10; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
